$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.09781768870467404
$ws.Range("C2").Value = 0.09781807909858378
$ws.Range("D2").Value = 0.220070122506433
$ws.Range("E2").Value = 0.4691163208698169
$ws.Range("F2").Value = 0.4691162829535249
$ws.Range("G2").Value = 23

$ws.Range("B3").Value = 0.3533564094378004
$ws.Range("C3").Value = 2.558293351813334
$ws.Range("D3").Value = 25.90499531260007
$ws.Range("E3").Value = 5.089695011746781
$ws.Range("F3").Value = 5.196899076927668
$ws.Range("G3").Value = 22

$ws.Range("B4").Value = -0.9761932629592093
$ws.Range("C4").Value = 2.467692221656369
$ws.Range("D4").Value = 19.72326690521151
$ws.Range("E4").Value = 4.441088482029097
$ws.Range("F4").Value = 4.439462726444811
$ws.Range("G4").Value = 21

$ws.Range("B5").Value = -0.1839716899615556
$ws.Range("C5").Value = 1.396256495667985
$ws.Range("D5").Value = 6.312839331709046
$ws.Range("E5").Value = 2.512536433906789
$ws.Range("F5").Value = 2.570888388128188
$ws.Range("G5").Value = 20

$ws.Range("B6").Value = -0.06823382884553478
$ws.Range("C6").Value = 1.510467036597221
$ws.Range("D6").Value = 9.867250468964778
$ws.Range("E6").Value = 3.141217991315595
$ws.Range("F6").Value = 3.226533207103521
$ws.Range("G6").Value = 19

$ws.Range("B7").Value = -0.3598138023919001
$ws.Range("C7").Value = 1.91869847179625
$ws.Range("D7").Value = 13.14814849056302
$ws.Range("E7").Value = 3.626037574345172
$ws.Range("F7").Value = 3.712746607591391
$ws.Range("G7").Value = 18

$ws.Range("B8").Value = -0.301078016622277
$ws.Range("C8").Value = 1.917698695039582
$ws.Range("D8").Value = 11.05059529941992
$ws.Range("E8").Value = 3.324243568004595
$ws.Range("F8").Value = 3.412468906127151
$ws.Range("G8").Value = 17

$ws.Range("B9").Value = -0.06847129159857163
$ws.Range("C9").Value = 1.755501846325137
$ws.Range("D9").Value = 11.33439003196579
$ws.Range("E9").Value = 3.366658585595782
$ws.Range("F9").Value = 3.476350839861456
$ws.Range("G9").Value = 16

$ws.Range("B10").Value = -0.1263805358386493
$ws.Range("C10").Value = 1.874779812578261
$ws.Range("D10").Value = 11.8364881125208
$ws.Range("E10").Value = 3.440419758186608
$ws.Range("F10").Value = 3.558769260475041
$ws.Range("G10").Value = 15

$ws.Range("B11").Value = -0.1131887169793432
$ws.Range("C11").Value = 2.30903469545812
$ws.Range("D11").Value = 15.32355388434682
$ws.Range("E11").Value = 3.914531119348372
$ws.Range("F11").Value = 4.072669564538719
